$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion-rate summary text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$old = $ws1.Range("A1").Value()
$new = $old.Replace("1000 Bs = 14.37 = 59399.48 pesos", "1000 Bs = 13.99 = 57748.39 pesos")
$new = $new.Replace("59399.48 pesos = 14.32 = 978.77 Bs", "57748.39 pesos = 13.96 = 967.78 Bs")
$ws1.Range("A1").Value = $new

# --- Sheet "tasas": update the rate table values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 71.5
$ws2.Range("O10").Value = 4129.01
$ws2.Range("N12").Value = 4137
$ws2.Range("O12").Value = 69.33
